$p = $ppt.ActivePresentation
$m = $p | Get-Member
Write-Output $m
